$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the student vaccination record sheet
$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Vaccine_Dose"

# Widen the Vaccine_Dose column (closest reachable width to 14.08984375)
$ws.Columns.Item(5).ColumnWidth = 13.33

# Cursor ends up one cell to the right of the typed headers, as in the source file
$ws.Range("F1").Select()
